{"js": "// Update the 20x5 table of addition/subtraction problems to the new\n// set of answers, then update the date/title line above it.\nconst body = context.document.body;\n\n// --- 1. Date/title paragraph (first paragraph of the body) ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2024-08-29 Thursday\", Word.InsertLocation.replace);\n\n// --- 2. Table of arithmetic problems (row-major values, formatting kept) ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\ntable.values = [\n  [\"91+1=92\", \"21-3=18\", \"74-69=5\", \"64+2=66\", \"15+36=51\"],\n  [\"85-79=6\", \"28-28=0\", \"21+74=95\", \"7-2=5\", \"3+91=94\"],\n  [\"97-8=89\", \"19+43=62\", \"95-40=55\", \"26+2=28\", \"74-54=20\"],\n  [\"71+15=86\", \"98-82=16\", \"93-49=44\", \"13+9=22\", \"31+28=59\"],\n  [\"55+9=64\", \"46+10=56\", \"11+25=36\", \"5+4=9\", \"38+3=41\"],\n  [\"51+11=62\", \"93-27=66\", \"17-3=14\", \"55-44=11\", \"12+51=63\"],\n  [\"81-38=43\", \"78+5=83\", \"75-5=70\", \"6+14=20\", \"5+26=31\"],\n  [\"20+50=70\", \"7+4=11\", \"37+58=95\", \"7+45=52\", \"91+1=92\"],\n  [\"35+42=77\", \"25+60=85\", \"44-21=23\", \"45-14=31\", \"99-30=69\"],\n  [\"70-13=57\", \"37+43=80\", \"28+7=35\", \"71-48=23\", \"0+25=25\"],\n  [\"86-57=29\", \"31+19=50\", \"91-67=24\", \"47+12=59\", \"91-12=79\"],\n  [\"18+39=57\", \"61-18=43\", \"94-29=65\", \"45-9=36\", \"36+13=49\"],\n  [\"7+87=94\", \"59-12=47\", \"56-19=37\", \"71-43=28\", \"97-1=96\"],\n  [\"64-25=39\", \"96+1=97\", \"21-14=7\", \"18+48=66\", \"94-76=18\"],\n  [\"80+1=81\", \"22+10=32\", \"52+38=90\", \"32+53=85\", \"73+10=83\"],\n  [\"29-25=4\", \"52+18=70\", \"93-66=27\", \"87-36=51\", \"85-2=83\"],\n  [\"86+6=92\", \"54-0=54\", \"72+21=93\", \"78-77=1\", \"7+69=76\"],\n  [\"2+97=99\", \"30+33=63\", \"98-93=5\", \"61+37=98\", \"6+57=63\"],\n  [\"2+56=58\", \"89+8=97\", \"56+12=68\", \"38+25=63\", \"80-12=68\"],\n  [\"26+60=86\", \"44-19=25\", \"6+37=43\", \"13+86=99\", \"84-9=75\"]\n];\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date/title paragraph (first paragraph of the body).\n$d.Paragraphs.Item(1).Range.Text = \"2024-08-29 Thursday\"\n\n# Update every arithmetic answer cell in the table, preserving formatting.\n$t = $d.Tables.Item(1)\n$answers = @(\n    @(\"91+1=92\", \"21-3=18\", \"74-69=5\", \"64+2=66\", \"15+36=51\"),\n    @(\"85-79=6\", \"28-28=0\", \"21+74=95\", \"7-2=5\", \"3+91=94\"),\n    @(\"97-8=89\", \"19+43=62\", \"95-40=55\", \"26+2=28\", \"74-54=20\"),\n    @(\"71+15=86\", \"98-82=16\", \"93-49=44\", \"13+9=22\", \"31+28=59\"),\n    @(\"55+9=64\", \"46+10=56\", \"11+25=36\", \"5+4=9\", \"38+3=41\"),\n    @(\"51+11=62\", \"93-27=66\", \"17-3=14\", \"55-44=11\", \"12+51=63\"),\n    @(\"81-38=43\", \"78+5=83\", \"75-5=70\", \"6+14=20\", \"5+26=31\"),\n    @(\"20+50=70\", \"7+4=11\", \"37+58=95\", \"7+45=52\", \"91+1=92\"),\n    @(\"35+42=77\", \"25+60=85\", \"44-21=23\", \"45-14=31\", \"99-30=69\"),\n    @(\"70-13=57\", \"37+43=80\", \"28+7=35\", \"71-48=23\", \"0+25=25\"),\n    @(\"86-57=29\", \"31+19=50\", \"91-67=24\", \"47+12=59\", \"91-12=79\"),\n    @(\"18+39=57\", \"61-18=43\", \"94-29=65\", \"45-9=36\", \"36+13=49\"),\n    @(\"7+87=94\", \"59-12=47\", \"56-19=37\", \"71-43=28\", \"97-1=96\"),\n    @(\"64-25=39\", \"96+1=97\", \"21-14=7\", \"18+48=66\", \"94-76=18\"),\n    @(\"80+1=81\", \"22+10=32\", \"52+38=90\", \"32+53=85\", \"73+10=83\"),\n    @(\"29-25=4\", \"52+18=70\", \"93-66=27\", \"87-36=51\", \"85-2=83\"),\n    @(\"86+6=92\", \"54-0=54\", \"72+21=93\", \"78-77=1\", \"7+69=76\"),\n    @(\"2+97=99\", \"30+33=63\", \"98-93=5\", \"61+37=98\", \"6+57=63\"),\n    @(\"2+56=58\", \"89+8=97\", \"56+12=68\", \"38+25=63\", \"80-12=68\"),\n    @(\"26+60=86\", \"44-19=25\", \"6+37=43\", \"13+86=99\", \"84-9=75\")\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $t.Cell($r, $c).Range.Text = $answers[$r - 1][$c - 1]\n    }\n}\n"}
